$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4595.5557
$ws.Range("I62").Value = 4979.5713
$ws.Range("J62").Value = 3251.5
$ws.Range("K62").Value = 4979.5713
$ws.Range("L62").Value = 3251.5
$ws.Range("M62").Value = -4355.5713
$ws.Range("N62").Value = -4499.5
$ws.Range("H65").Value = 4595.5557
$ws.Range("I65").Value = 4979.5713
$ws.Range("J65").Value = 3251.5
$ws.Range("K65").Value = 24897.8565
$ws.Range("L65").Value = 16257.5
$ws.Range("M65").Value = -21777.8565
$ws.Range("N65").Value = -22497.5
$ws.Range("H69").Value = 3506
$ws.Range("I69").Value = 3013
$ws.Range("K69").Value = 9039
$ws.Range("M69").Value = -8165
$ws.Range("H72").Value = 3506
$ws.Range("I72").Value = 3013
$ws.Range("K72").Value = 27117
$ws.Range("M72").Value = -22749
$ws.Range("H100").Value = 1873.5518
$ws.Range("I100").Value = 1569.2106
$ws.Range("K100").Value = 1569.2106
$ws.Range("M100").Value = -1028.2106
$ws.Range("H141").Value = 4493.3887
$ws.Range("I141").Value = 2897.5833
$ws.Range("J141").Value = 7685
$ws.Range("K141").Value = 8692.749899999999
$ws.Range("L141").Value = 23055
$ws.Range("M141").Value = -3512.749899999999
$ws.Range("N141").Value = -33415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3410.3333
$ws.Range("I63").Value = 2308.3333
$ws.Range("J63").Value = 5063.3335
$ws.Range("K63").Value = 2308.3333
$ws.Range("L63").Value = 5063.3335
$ws.Range("M63").Value = -1622.3333
$ws.Range("N63").Value = -6435.3335
$ws.Range("H66").Value = 3410.3333
$ws.Range("I66").Value = 2308.3333
$ws.Range("J66").Value = 5063.3335
$ws.Range("K66").Value = 11541.6665
$ws.Range("L66").Value = 25316.6675
$ws.Range("M66").Value = -8109.666499999999
$ws.Range("N66").Value = -32180.6675
$ws.Range("H74").Value = 2011.7046
$ws.Range("I74").Value = 1712.7428
$ws.Range("K74").Value = 1712.7428
$ws.Range("M74").Value = -838.7428
$ws.Range("H77").Value = 2011.7046
$ws.Range("I77").Value = 1712.7428
$ws.Range("K77").Value = 8563.714
$ws.Range("M77").Value = -4195.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3598.7673
$ws.Range("J134").Value = 3432.4856
$ws.Range("L134").Value = 10297.4568
$ws.Range("N134").Value = -15367.4568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 26930.045
$ws.Range("J26").Value = 26930.045
$ws.Range("L26").Value = 26930.045
$ws.Range("N26").Value = -27504.045
$ws.Range("H31").Value = 5381521.5
$ws.Range("I31").Value = 2400.476
$ws.Range("K31").Value = 2400.476
$ws.Range("M31").Value = -2105.476
$ws.Range("H34").Value = 5381521.5
$ws.Range("I34").Value = 2400.476
$ws.Range("K34").Value = 2400.476
$ws.Range("M34").Value = -2198.476
$ws.Range("H132").Value = 46238.906
$ws.Range("I132").Value = 1287.4286
$ws.Range("J132").Value = 132055.36
$ws.Range("K132").Value = 3862.2858
$ws.Range("L132").Value = 396166.08
$ws.Range("M132").Value = -1332.2858
$ws.Range("N132").Value = -401226.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4023.889
$ws.Range("J75").Value = 4023.889
$ws.Range("L75").Value = 12071.667
$ws.Range("N75").Value = -14067.667
$ws.Range("H78").Value = 4023.889
$ws.Range("J78").Value = 4023.889
$ws.Range("L78").Value = 36215.001
$ws.Range("N78").Value = -46199.001
$ws.Range("H97").Value = 18998.045
$ws.Range("I97").Value = 491.33334
$ws.Range("J97").Value = 25938.062
$ws.Range("K97").Value = 1474.00002
$ws.Range("L97").Value = 77814.186
$ws.Range("M97").Value = -978.0000199999999
$ws.Range("N97").Value = -78806.186
$ws.Range("H98").Value = 618.25
$ws.Range("J98").Value = 638.0909
$ws.Range("L98").Value = 1914.2727
$ws.Range("N98").Value = -4910.2727
$ws.Range("H101").Value = 6000
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H107").Value = 508.66666
$ws.Range("I107").Value = 513.8333
$ws.Range("J107").Value = 503.5
$ws.Range("K107").Value = 1541.4999
$ws.Range("L107").Value = 1510.5
$ws.Range("M107").Value = 378.5001
$ws.Range("N107").Value = -5350.5
$ws.Range("H110").Value = 2506.75
$ws.Range("I110").Value = 1263.5
$ws.Range("J110").Value = 3750
$ws.Range("K110").Value = 3790.5
$ws.Range("L110").Value = 11250
$ws.Range("M110").Value = 299.5
$ws.Range("N110").Value = -19430
$ws.Range("H131").Value = 852.6598
$ws.Range("I131").Value = 503.33334
$ws.Range("J131").Value = 863.80853
$ws.Range("K131").Value = 1510.00002
$ws.Range("L131").Value = 2591.42559
$ws.Range("M131").Value = 3529.99998
$ws.Range("N131").Value = -12671.42559

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7023.75
$ws.Range("I80").Value = 5450
$ws.Range("J80").Value = 8597.5
$ws.Range("K80").Value = 5450
$ws.Range("L80").Value = 8597.5
$ws.Range("M80").Value = -4452
$ws.Range("N80").Value = -10593.5
$ws.Range("H83").Value = 7023.75
$ws.Range("I83").Value = 5450
$ws.Range("J83").Value = 8597.5
$ws.Range("K83").Value = 27250
$ws.Range("L83").Value = 42987.5
$ws.Range("M83").Value = -22258
$ws.Range("N83").Value = -52971.5
$ws.Range("H132").Value = 83337140
$ws.Range("I132").Value = 200002190
$ws.Range("J132").Value = 4943.2856
$ws.Range("K132").Value = 600006570
$ws.Range("L132").Value = 14829.8568
$ws.Range("M132").Value = -600004040
$ws.Range("N132").Value = -19889.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 27778946
$ws.Range("I82").Value = 1750
$ws.Range("J82").Value = 83333336
$ws.Range("K82").Value = 1750
$ws.Range("L82").Value = 83333336
$ws.Range("M82").Value = -1389
$ws.Range("N82").Value = -83334058
$ws.Range("H85").Value = 27778946
$ws.Range("I85").Value = 1750
$ws.Range("J85").Value = 83333336
$ws.Range("K85").Value = 1750
$ws.Range("L85").Value = 83333336
$ws.Range("M85").Value = -502
$ws.Range("N85").Value = -83335832
$ws.Range("H93").Value = 1003.2105
$ws.Range("I93").Value = 926.1
$ws.Range("J93").Value = 1088.8889
$ws.Range("K93").Value = 926.1
$ws.Range("L93").Value = 1088.8889
$ws.Range("M93").Value = 321.9
$ws.Range("N93").Value = -3584.8889
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 145648.14
$ws.Range("I122").Value = 252784.25
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 758352.75
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -755902.75
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1178.7778
$ws.Range("I81").Value = 1034.8334
$ws.Range("J81").Value = 1466.6666
$ws.Range("K81").Value = 2069.6668
$ws.Range("L81").Value = 2933.3332
$ws.Range("M81").Value = -1008.6668
$ws.Range("N81").Value = -5055.3332
$ws.Range("H84").Value = 1178.7778
$ws.Range("I84").Value = 1034.8334
$ws.Range("J84").Value = 1466.6666
$ws.Range("K84").Value = 10348.334
$ws.Range("L84").Value = 14666.666
$ws.Range("M84").Value = -5044.333999999999
$ws.Range("N84").Value = -25274.666
$ws.Range("H119").Value = 42146.5
$ws.Range("J119").Value = 42146.5
$ws.Range("L119").Value = 42146.5
$ws.Range("N119").Value = -51822.5
$ws.Range("H122").Value = 2381869.5
$ws.Range("I122").Value = 2598312
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7794936
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -7792486
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 1177588.8
$ws.Range("I126").Value = 1548657.8
$ws.Range("K126").Value = 4645973.4
$ws.Range("M126").Value = -4643503.4
